$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) / Volume(1h) (E) columns for existing rows 2-42 ---
# Row 2
$ws.Range("D2").Value = '69.357.94'
$ws.Range("E2").Value = '  +1.70%  '
# Row 3
$ws.Range("D3").Value = '3.378.20'
$ws.Range("E3").Value = '  +1.22%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.01%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.49'
$ws.Range("E5").Value = '  -0.40%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.58'
$ws.Range("E6").Value = '  +0.83%  '
# Row 7
$ws.Range("E7").Value = '  +0.05%  '
# Row 8
$ws.Range("E8").Value = '  +0.56%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.198'
$ws.Range("E9").Value = '  +8.57%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.587'
$ws.Range("E10").Value = '  +0.83%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.43'
$ws.Range("E11").Value = '  +0.88%  '
# Row 12
$ws.Range("E12").Value = '  +3.95%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '687.93'
$ws.Range("E13").Value = '  -0.84%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.60'
$ws.Range("E14").Value = '  +2.19%  '
# Row 15
$ws.Range("D15").Value = '3.918.28'
$ws.Range("E15").Value = '  +0.87%  '
# Row 16
$ws.Range("D16").Value = '69.432.00'
$ws.Range("E16").Value = '  +1.77%  '
# Row 17
$ws.Range("E17").Value = '  +0.80%  '
# Row 18
$ws.Range("D18").Value = '3.382.27'
$ws.Range("E18").Value = '  +1.28%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.78'
$ws.Range("E19").Value = '  +2.00%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.26'
$ws.Range("E20").Value = '  +1.02%  '
# Row 21
$ws.Range("E21").Value = '  +1.59%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.20'
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.35'
$ws.Range("E23").Value = '  -2.02%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '101.44'
$ws.Range("E24").Value = '  +1.42%  '
# Row 25
$ws.Range("E25").Value = '  -0.69%  '
# Row 26
$ws.Range("E26").Value = '  -0.06%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.70'
$ws.Range("E27").Value = '  +2.08%  '
# Row 28
$ws.Range("E28").Value = '  +1.68%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.72'
$ws.Range("E29").Value = '  +2.75%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.91'
$ws.Range("E30").Value = '  -0.20%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.82'
$ws.Range("E31").Value = '  +16.95%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.04'
$ws.Range("E32").Value = '  +0.18%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '550.46'
$ws.Range("E33").Value = '  -2.33%  '
# Row 34
$ws.Range("E34").Value = '  +0.26%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.90'
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.05%  '
# Row 37
$ws.Range("D37").Value = '3.603.70'
$ws.Range("E37").Value = '  -2.36%  '
# Row 38
$ws.Range("E38").Value = '  +3.32%  '
# Row 39
$ws.Range("E39").Value = '  +1.49%  '
# Row 40
$ws.Range("D40").Value = '0.0₃0727'
$ws.Range("E40").Value = '  +8.42%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.31'
$ws.Range("E41").Value = '  +4.61%  '
# Row 42
$ws.Range("E42").Value = '  +4.35%  '

# --- Rows 43-51: a new coin (ApeXProtocol) was inserted at row 43, shifting the
# remaining coins down by one; THORChain (formerly row 51) drops off the list.
# Column A (the numeric index) is untouched, so only B/C/D/E are rewritten directly. ---
# Row 43: ApeXProtocol
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.38'
$ws.Range("E43").Value = '  +3.92%  '
# Row 44: VeChain
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0424'
$ws.Range("E44").Value = '  +2.76%  '
# Row 45: TheGraph
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.336'
$ws.Range("E45").Value = '  +0.36%  '
# Row 46: ThetaToken
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.65'
$ws.Range("E46").Value = '  +0.28%  '
# Row 47: Stellar
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.129'
$ws.Range("E47").Value = '  +0.41%  '
# Row 48: FirstDigitalUSD
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").Value = '  -0.24%  '
# Row 49: Mantle
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.38'
$ws.Range("E49").Value = '  +3.68%  '
# Row 50: Monero
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '129.22'
$ws.Range("E50").Value = '  -1.29%  '
# Row 51: CoreDAO
$ws.Range("B51").Value = 'CoreDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.58'
$ws.Range("E51").Value = '  +0.35%  '
